# Update TPM-derived values in the LR-pairs sheet to reflect the new TPM
# computations (commit: "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 0.07271233333333334
$ws.Range("N2").Value = 0.218137
$ws.Range("O2").Value = 0.004171225362010892
$ws.Range("P2").Value = 0.004171225362010893
$ws.Range("Q2").Value = 0.02874585148555556
$ws.Range("R2").Value = 0.25871266337
$ws.Range("S2").Value = 0.002308522384773037
$ws.Range("T2").Value = 0.002308522384773038

# Row 3
$ws.Range("O3").Value = 0.5387060579248023
$ws.Range("P3").Value = 0.5387060579248023
$ws.Range("S3").Value = 0.2981414058464384
$ws.Range("T3").Value = 0.2981414058464384

# Row 4
$ws.Range("O4").Value = 0.4571227167131868
$ws.Range("P4").Value = 0.4571227167131868
$ws.Range("Q4").Value = 3.150244972373334
$ws.Range("S4").Value = 0.2529899328220233
$ws.Range("T4").Value = 0.2529899328220233

# Row 5
$ws.Range("M5").Value = 0.07271233333333334
$ws.Range("N5").Value = 0.218137
$ws.Range("O5").Value = 0.004171225362010892
$ws.Range("P5").Value = 0.004171225362010893
$ws.Range("Q5").Value = 0.02319448297255556
$ws.Range("R5").Value = 0.208750346753
$ws.Range("S5").Value = 0.001862702977237855
$ws.Range("T5").Value = 0.001862702977237855

# Row 6
$ws.Range("O6").Value = 0.5387060579248023
$ws.Range("P6").Value = 0.5387060579248023
$ws.Range("S6").Value = 0.2405646520783639
$ws.Range("T6").Value = 0.2405646520783639

# Row 7
$ws.Range("O7").Value = 0.4571227167131868
$ws.Range("P7").Value = 0.4571227167131868
$ws.Range("S7").Value = 0.2041327838911635
$ws.Range("T7").Value = 0.2041327838911635

$wb.Save()
